# Auto-generated edit script: apply value updates to Leviathan_Profits workbook
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H88").Value = 84825.086
$ws.Range("I88").Value = 1566
$ws.Range("J88").Value = 112578.11
$ws.Range("K88").Value = 1566
$ws.Range("L88").Value = 112578.11
$ws.Range("M88").Value = -1160
$ws.Range("N88").Value = -113390.11
$ws.Range("H91").Value = 84825.086
$ws.Range("I91").Value = 1566
$ws.Range("J91").Value = 112578.11
$ws.Range("K91").Value = 1566
$ws.Range("L91").Value = 112578.11
$ws.Range("M91").Value = -162
$ws.Range("N91").Value = -115386.11
$ws.Range("H103").Value = 590.9091
$ws.Range("J103").Value = 590.9091
$ws.Range("L103").Value = 1772.7273
$ws.Range("N103").Value = -2944.7273
$ws.Range("H129").Value = 1128.0625
$ws.Range("I129").Value = 665.4
$ws.Range("J129").Value = 1899.1666
$ws.Range("K129").Value = 1996.2
$ws.Range("L129").Value = 5697.4998
$ws.Range("M129").Value = 3003.8
$ws.Range("N129").Value = -15697.4998
$ws.Range("H132").Value = 3661.162
$ws.Range("I132").Value = 1263.72
$ws.Range("J132").Value = 8655.833000000001
$ws.Range("K132").Value = 3791.16
$ws.Range("L132").Value = 25967.499
$ws.Range("M132").Value = -1261.16
$ws.Range("N132").Value = -31027.499
$ws.Range("H137").Value = 60646
$ws.Range("I137").Value = 1883.4615
$ws.Range("J137").Value = 251624.25
$ws.Range("K137").Value = 5650.3845
$ws.Range("L137").Value = 754872.75
$ws.Range("M137").Value = -3100.3845
$ws.Range("N137").Value = -759972.75
$ws.Range("H138").Value = 2136.4443
$ws.Range("I138").Value = 1639.25
$ws.Range("J138").Value = 3130.8333
$ws.Range("K138").Value = 4917.75
$ws.Range("L138").Value = 9392.499899999999
$ws.Range("M138").Value = 222.25
$ws.Range("N138").Value = -19672.4999
$ws.Range("H141").Value = 77566
$ws.Range("I141").Value = 102375.78
$ws.Range("K141").Value = 307127.34
$ws.Range("M141").Value = -301947.34
$ws.Range("N141").ClearContents()

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 716.5
$ws.Range("I110").Value = 716.5
$ws.Range("J110").Value = 0
$ws.Range("K110").Value = 716.5
$ws.Range("L110").Value = 0
$ws.Range("M110").Value = 1328.5
$ws.Range("N110").ClearContents()
$ws.Range("H132").Value = 2899.08
$ws.Range("I132").Value = 1991.6364
$ws.Range("K132").Value = 5974.9092
$ws.Range("M132").Value = -3444.9092
$ws.Range("N132").ClearContents()

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 538.8421
$ws.Range("I94").Value = 544.625
$ws.Range("K94").Value = 544.625
$ws.Range("M94").Value = -93.625
$ws.Range("N94").ClearContents()
$ws.Range("H105").Value = 3706090.2
$ws.Range("I105").Value = 6252151
$ws.Range("J105").Value = 2729
$ws.Range("K105").Value = 6252151
$ws.Range("L105").Value = 2729
$ws.Range("M105").Value = -6250404
$ws.Range("N105").Value = -6223
$ws.Range("H134").Value = 3236.0715
$ws.Range("I134").Value = 2876.25
$ws.Range("K134").Value = 8628.75
$ws.Range("M134").Value = -6093.75
$ws.Range("N134").ClearContents()

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1661.8572
$ws.Range("I31").Value = 1661.8572
$ws.Range("K31").Value = 1661.8572
$ws.Range("M31").Value = -1366.8572
$ws.Range("H34").Value = 1661.8572
$ws.Range("I34").Value = 1661.8572
$ws.Range("K34").Value = 1661.8572
$ws.Range("M34").Value = -1459.8572
$ws.Range("H51").Value = 9999.200000000001
$ws.Range("J51").Value = 9999.200000000001
$ws.Range("L51").Value = 9999.200000000001
$ws.Range("N51").Value = -11471.2
$ws.Range("H58").Value = 1934.6364
$ws.Range("I58").Value = 1630.3334
$ws.Range("K58").Value = 1630.3334
$ws.Range("M58").Value = -1427.3334
$ws.Range("N58").ClearContents()
$ws.Range("H61").Value = 9999.200000000001
$ws.Range("J61").Value = 9999.200000000001
$ws.Range("L61").Value = 9999.200000000001
$ws.Range("N61").Value = -10695.2
$ws.Range("H74").Value = 22222
$ws.Range("J74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("N74").ClearContents()
$ws.Range("H77").Value = 22222
$ws.Range("J77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("N77").ClearContents()
$ws.Range("H132").Value = 1907.9117
$ws.Range("I132").Value = 1899.9678
$ws.Range("J132").Value = 1990
$ws.Range("K132").Value = 5699.903399999999
$ws.Range("L132").Value = 5970
$ws.Range("M132").Value = -3169.903399999999
$ws.Range("N132").Value = -11030
$ws.Range("H134").Value = 1945.303
$ws.Range("I134").Value = 1828.7037
$ws.Range("K134").Value = 5486.1111
$ws.Range("M134").Value = -2951.1111
$ws.Range("N134").ClearContents()
$ws.Range("H136").Value = 1934.6364
$ws.Range("I136").Value = 1630.3334
$ws.Range("K136").Value = 4891.0002
$ws.Range("M136").Value = -2341.0002
$ws.Range("N136").ClearContents()

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 360.75
$ws.Range("I14").Value = 360.75
$ws.Range("K14").Value = 1082.25
$ws.Range("M14").Value = -909.25
$ws.Range("H23").Value = 0
$ws.Range("J23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("N23").ClearContents()
$ws.Range("H32").Value = 928175.7
$ws.Range("J32").Value = 1853851.4
$ws.Range("L32").Value = 5561554.199999999
$ws.Range("N32").Value = -5562120.199999999
$ws.Range("H109").Value = 27
$ws.Range("I109").Value = 27
$ws.Range("K109").Value = 81
$ws.Range("M109").Value = 959
$ws.Range("H122").Value = 789.36365
$ws.Range("I122").Value = 195
$ws.Range("K122").Value = 1755
$ws.Range("M122").Value = 695
$ws.Range("H133").Value = 10665.167
$ws.Range("I133").Value = 6996.6
$ws.Range("J133").Value = 12076.154
$ws.Range("K133").Value = 20989.8
$ws.Range("L133").Value = 36228.462
$ws.Range("M133").Value = -15929.8
$ws.Range("N133").Value = -46348.462

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H17").Value = 2955.2727
$ws.Range("J17").Value = 3150
$ws.Range("L17").Value = 3150
$ws.Range("N17").Value = -3486
$ws.Range("H122").Value = 1484.8
$ws.Range("I122").Value = 1349.9286
$ws.Range("J122").Value = 1799.5
$ws.Range("K122").Value = 4049.7858
$ws.Range("L122").Value = 5398.5
$ws.Range("M122").Value = -1599.7858
$ws.Range("N122").Value = -10298.5
$ws.Range("H126").Value = 3131.375
$ws.Range("I126").Value = 2249
$ws.Range("K126").Value = 6747
$ws.Range("M126").Value = -4277
$ws.Range("N126").ClearContents()
$ws.Range("H136").Value = 52775
$ws.Range("J136").Value = 52775
$ws.Range("L136").Value = 158325
$ws.Range("N136").Value = -163425

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1585.8889
$ws.Range("I22").Value = 1609.4
$ws.Range("J22").Value = 1556.5
$ws.Range("K22").Value = 1609.4
$ws.Range("L22").Value = 1556.5
$ws.Range("M22").Value = -1314.4
$ws.Range("N22").Value = -2146.5
$ws.Range("H27").Value = 1585.8889
$ws.Range("I27").Value = 1609.4
$ws.Range("J27").Value = 1556.5
$ws.Range("K27").Value = 1609.4
$ws.Range("L27").Value = 1556.5
$ws.Range("M27").Value = -1502.4
$ws.Range("N27").Value = -1770.5
$ws.Range("H55").Value = 594.9091
$ws.Range("I55").Value = 520
$ws.Range("J55").Value = 623
$ws.Range("K55").Value = 520
$ws.Range("L55").Value = 623
$ws.Range("M55").Value = -347
$ws.Range("N55").Value = -969
$ws.Range("H61").Value = 9824
$ws.Range("I61").Value = 12643.637
$ws.Range("J61").Value = 3620.8
$ws.Range("K61").Value = 12643.637
$ws.Range("L61").Value = 3620.8
$ws.Range("M61").Value = -12441.637
$ws.Range("N61").Value = -4024.8
$ws.Range("H68").Value = 3399.3333
$ws.Range("I68").Value = 3399.3333
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 3399.3333
$ws.Range("L68").Value = 0
$ws.Range("M68").Value = -2650.3333
$ws.Range("N68").ClearContents()
$ws.Range("H71").Value = 3399.3333
$ws.Range("I71").Value = 3399.3333
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 16996.6665
$ws.Range("L71").Value = 0
$ws.Range("M71").Value = -13252.6665
$ws.Range("N71").ClearContents()
$ws.Range("H93").Value = 1793
$ws.Range("I93").Value = 1672.2
$ws.Range("K93").Value = 1672.2
$ws.Range("M93").Value = -424.2
$ws.Range("N93").ClearContents()
$ws.Range("H113").Value = 9824
$ws.Range("I113").Value = 12643.637
$ws.Range("J113").Value = 3620.8
$ws.Range("K113").Value = 12643.637
$ws.Range("L113").Value = 3620.8
$ws.Range("M113").Value = -10473.637
$ws.Range("N113").Value = -7960.8
$ws.Range("H132").Value = 3035.7942
$ws.Range("I132").Value = 2915.8572
$ws.Range("J132").Value = 3229.5386
$ws.Range("K132").Value = 8747.571599999999
$ws.Range("L132").Value = 9688.6158
$ws.Range("M132").Value = -6217.571599999999
$ws.Range("N132").Value = -14748.6158

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 256374.75
$ws.Range("I62").Value = 8250
$ws.Range("K62").Value = 8250
$ws.Range("M62").Value = -7626
$ws.Range("N62").ClearContents()
$ws.Range("H65").Value = 256374.75
$ws.Range("I65").Value = 8250
$ws.Range("K65").Value = 41250
$ws.Range("M65").Value = -38130
$ws.Range("N65").ClearContents()
$ws.Range("H132").Value = 12127.789
$ws.Range("I132").Value = 14729.929
$ws.Range("K132").Value = 44189.787
$ws.Range("M132").Value = -41659.787
$ws.Range("N132").ClearContents()
